# Quarterly indexing esoteric bug-fix operation
#
# Column A holds the period-end date stamps for the nowcast series. They were
# previously stored as "end of period, minus one second" serials (e.g.
# 40543.99999999999 == 2010-12-31 23:59:59.999999). The fix re-expresses each
# stamp as the 15th of the appropriate month (the true quarterly/semiannual
# index date), stored as a clean integer day-serial with no fractional part.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newSerials = @{
    2  = 40497
    3  = 40678
    4  = 40862
    5  = 41044
    6  = 41228
    7  = 41409
    8  = 41593
    9  = 41774
    10 = 41958
    11 = 42139
    12 = 42323
    13 = 42505
    14 = 42689
    15 = 42870
    16 = 43054
    17 = 43146
    18 = 43235
    19 = 43327
    20 = 43419
    21 = 43511
    22 = 43600
    23 = 43692
    24 = 43784
    25 = 43876
    26 = 43966
    27 = 44058
    28 = 44150
    29 = 44242
    30 = 44331
    31 = 44423
    32 = 44515
    33 = 44607
    34 = 44696
    35 = 44788
    36 = 44880
    37 = 44972
    38 = 45061
    39 = 45153
    40 = 45245
    41 = 45337
    42 = 45427
    43 = 45519
    44 = 45611
    45 = 45703
    46 = 45792
    47 = 45884
}

foreach ($row in $newSerials.Keys) {
    $ws.Cells.Item($row, 1).Value = $newSerials[$row]
}
